# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1361
$ws1.Range("F7").Value  = 11689
$ws1.Range("F8").Value  = 4388
$ws1.Range("F10").Value = 39
$ws1.Range("F13").Value = 2547
$ws1.Range("F15").Value = 147
$ws1.Range("F16").Value = 40
$ws1.Range("F17").Value = 5099
$ws1.Range("F20").Value = 515
$ws1.Range("F22").Value = 11273

# Sheet "全部类型" (sheet4): same events, rows shifted by one vs "展览"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1361
$ws4.Range("F7").Value  = 11689
$ws4.Range("F8").Value  = 4388
$ws4.Range("F10").Value = 39
$ws4.Range("F13").Value = 2547
$ws4.Range("F16").Value = 147
$ws4.Range("F17").Value = 40
$ws4.Range("F18").Value = 5099
$ws4.Range("F21").Value = 515
$ws4.Range("F23").Value = 11273
